# Regenerate s_val data to filter save games.
# Updates the numeric stat columns (B:G) for rows 2-6 on the active sheet
# with the newly computed values, leaving column A (dates) and F (Win flag)
# untouched where they did not change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1169995834814548, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.426980108624251)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    4 = @(1.445647641019636, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 0, 5.507293877332936)
    5 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    6 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F - Win
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G - sum
}
